$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "remaining days" (column E) simply ticks down by one day,
# reflecting the daily auto-update of the tracker.
$decrementUpdates = @{
    2 = 6
    3 = 6
    4 = 6
    5 = 2
    6 = 6
    7 = 2
    8 = 6
    9 = 2
    10 = 6
    11 = 6
    12 = 2
    13 = 6
    14 = 6
    15 = 6
    16 = 6
    17 = 2
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 2
    23 = 2
    24 = 2
    25 = 2
    26 = 2
    28 = 5
    29 = 5
    30 = 5
    31 = 5
    32 = 5
    33 = 5
    34 = 5
    35 = 5
    37 = 5
    38 = 5
    39 = 5
    40 = 6
    41 = 6
    42 = 5
    43 = 2
    44 = 6
    45 = 2
    46 = 6
    47 = 5
    48 = 6
    58 = 4
    59 = 4
    60 = 4
    62 = 4
    63 = 4
    64 = 4
    65 = 5
    66 = 5
    67 = 5
    68 = 5
    69 = 5
    70 = 6
    71 = 6
    72 = 6
    73 = 6
    74 = 6
    75 = 6
    76 = 6
    77 = 9
    78 = 9
    79 = 9
    80 = 9
    81 = 9
    82 = 9
    83 = 9
    84 = 9
    85 = 9
    86 = 9
    87 = 6
    88 = 6
    89 = 6
    90 = 6
    91 = 2
    92 = 6
    93 = 9
    94 = 2
    95 = 8
    96 = 6
    97 = 6
    98 = 6
    99 = 6
}

foreach ($row in $decrementUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $decrementUpdates[$row]
}

# Rows where the item was restocked: "remaining days" (E) resets back up
# and the "start date" (F) is set to the restock date 2025-12-23.
$resetUpdates = @(
    @{ Row = 27; E = 7; F = 20251223 }
    @{ Row = 49; E = 7; F = 20251223 }
    @{ Row = 50; E = 10; F = 20251223 }
    @{ Row = 51; E = 10; F = 20251223 }
    @{ Row = 52; E = 10; F = 20251223 }
    @{ Row = 53; E = 10; F = 20251223 }
    @{ Row = 54; E = 10; F = 20251223 }
    @{ Row = 55; E = 10; F = 20251223 }
    @{ Row = 56; E = 10; F = 20251223 }
    @{ Row = 57; E = 10; F = 20251223 }
    @{ Row = 61; E = 7; F = 20251223 }
)

foreach ($item in $resetUpdates) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    $ws.Cells.Item($item.Row, 6).Value = $item.F
}

$wb.Save()